{"js": "// Update the date/title line and every arithmetic expression in the\n// practice-sheet table, preserving each run's existing formatting by\n// using paragraph-level insertText(..., Word.InsertLocation.replace)\n// instead of wiping/recreating runs.\n\nconst newTitle = \"2023-12-28 Thursday\";\n\n// New cell text, row-major (20 rows x 5 columns), matching the table's\n// existing reading order so duplicate old values (e.g. \"93-17=\" appears\n// twice) are disambiguated correctly by position rather than by text.\nconst newGrid = [\n  [\"77-72=\", \"16+41=\", \"29+34=\", \"46+23=\", \"82-18=\"],\n  [\"96-40=\", \"23+40=\", \"55+42=\", \"25+51=\", \"72-70=\"],\n  [\"45+5=\", \"56+7=\", \"74-36=\", \"39-39=\", \"27+19=\"],\n  [\"39+18=\", \"64-12=\", \"27+64=\", \"88+7=\", \"40+48=\"],\n  [\"57+35=\", \"44+50=\", \"91-10=\", \"89-88=\", \"26-26=\"],\n  [\"60+14=\", \"27-4=\", \"70-22=\", \"88-23=\", \"54+20=\"],\n  [\"94-17=\", \"63-23=\", \"5+57=\", \"23+19=\", \"3+77=\"],\n  [\"90+7=\", \"57-20=\", \"89-21=\", \"41-20=\", \"39-25=\"],\n  [\"32+42=\", \"76+12=\", \"43-35=\", \"15+58=\", \"42+15=\"],\n  [\"90-48=\", \"43-15=\", \"1+78=\", \"65-13=\", \"82-13=\"],\n  [\"94-43=\", \"49-32=\", \"31+25=\", \"68+4=\", \"13+79=\"],\n  [\"29-22=\", \"88+8=\", \"14+13=\", \"54-8=\", \"25-24=\"],\n  [\"16+76=\", \"17+73=\", \"63+34=\", \"74-0=\", \"52-29=\"],\n  [\"77+14=\", \"50-8=\", \"37-32=\", \"43+5=\", \"9+79=\"],\n  [\"97-90=\", \"6+11=\", \"68-10=\", \"88-55=\", \"53-45=\"],\n  [\"98-7=\", \"80-30=\", \"28-10=\", \"3+96=\", \"0+87=\"],\n  [\"5+21=\", \"65-65=\", \"57+0=\", \"34+64=\", \"86-82=\"],\n  [\"66-65=\", \"68+29=\", \"66-8=\", \"25+31=\", \"29+50=\"],\n  [\"94-36=\", \"15+57=\", \"98-79=\", \"58+26=\", \"95-5=\"],\n  [\"98-50=\", \"71-27=\", \"19+7=\", \"20+59=\", \"27+52=\"],\n];\n\nconst body = context.document.body;\n\n// 1) Title / date paragraph (first paragraph in the body, outside the\n//    table) \u2014 replace its text in place so its run formatting (Arial,\n//    sz 30) survives untouched.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nconst titleParagraph = paragraphs.items[0];\ntitleParagraph.insertText(newTitle, Word.InsertLocation.replace);\n\n// 2) The practice table: walk every row/column and replace the single\n//    paragraph's text inside each cell, again via insertText(replace)\n//    so the TimeNewRoman/sz-30 run formatting and the cell's left\n//    paragraph alignment are preserved.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nfor (let r = 0; r < newGrid.length; r++) {\n  for (let c = 0; c < newGrid[r].length; c++) {\n    const cell = table.getCell(r, c);\n    const cellParagraphs = cell.body.paragraphs;\n    cellParagraphs.load(\"items\");\n    await context.sync();\n    const cellParagraph = cellParagraphs.items[0];\n    cellParagraph.insertText(newGrid[r][c], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date/title line and every arithmetic expression in the\n# practice-sheet table, preserving each run's existing formatting by\n# assigning to Range.Text (which replaces only the content, leaving the\n# paragraph mark / cell mark and the run's rPr / paragraph's pPr intact).\n\n$d = $word.ActiveDocument\n\n# New title/date text.\n$newTitle = \"2023-12-28 Thursday\"\n\n# New cell text, row-major (20 rows x 5 columns), matching the table's\n# existing reading order so duplicate old values (e.g. \"93-17=\" appears\n# twice) are disambiguated correctly by position rather than by text.\n$newValues = @(\n    @(\"77-72=\", \"16+41=\", \"29+34=\", \"46+23=\", \"82-18=\"),\n    @(\"96-40=\", \"23+40=\", \"55+42=\", \"25+51=\", \"72-70=\"),\n    @(\"45+5=\", \"56+7=\", \"74-36=\", \"39-39=\", \"27+19=\"),\n    @(\"39+18=\", \"64-12=\", \"27+64=\", \"88+7=\", \"40+48=\"),\n    @(\"57+35=\", \"44+50=\", \"91-10=\", \"89-88=\", \"26-26=\"),\n    @(\"60+14=\", \"27-4=\", \"70-22=\", \"88-23=\", \"54+20=\"),\n    @(\"94-17=\", \"63-23=\", \"5+57=\", \"23+19=\", \"3+77=\"),\n    @(\"90+7=\", \"57-20=\", \"89-21=\", \"41-20=\", \"39-25=\"),\n    @(\"32+42=\", \"76+12=\", \"43-35=\", \"15+58=\", \"42+15=\"),\n    @(\"90-48=\", \"43-15=\", \"1+78=\", \"65-13=\", \"82-13=\"),\n    @(\"94-43=\", \"49-32=\", \"31+25=\", \"68+4=\", \"13+79=\"),\n    @(\"29-22=\", \"88+8=\", \"14+13=\", \"54-8=\", \"25-24=\"),\n    @(\"16+76=\", \"17+73=\", \"63+34=\", \"74-0=\", \"52-29=\"),\n    @(\"77+14=\", \"50-8=\", \"37-32=\", \"43+5=\", \"9+79=\"),\n    @(\"97-90=\", \"6+11=\", \"68-10=\", \"88-55=\", \"53-45=\"),\n    @(\"98-7=\", \"80-30=\", \"28-10=\", \"3+96=\", \"0+87=\"),\n    @(\"5+21=\", \"65-65=\", \"57+0=\", \"34+64=\", \"86-82=\"),\n    @(\"66-65=\", \"68+29=\", \"66-8=\", \"25+31=\", \"29+50=\"),\n    @(\"94-36=\", \"15+57=\", \"98-79=\", \"58+26=\", \"95-5=\"),\n    @(\"98-50=\", \"71-27=\", \"19+7=\", \"20+59=\", \"27+52=\")\n)\n\n# 1) Title / date paragraph (first paragraph in the body, outside the\n#    table) \u2014 assigning Range.Text replaces the run's text while the\n#    run's formatting (Arial, sz 30) and paragraph mark are untouched.\n$titleRange = $d.Paragraphs.Item(1).Range\n$titleRange.Text = $newTitle\n\n# 2) The practice table: walk every row/column and overwrite the cell's\n#    Range.Text, which preserves the TimeNewRoman/sz-30 run formatting\n#    and the cell's left paragraph alignment (and the cell end mark).\n$table = $d.Tables.Item(1)\n$rowCount = $newValues.Count\nfor ($r = 1; $r -le $rowCount; $r++) {\n    $rowValues = $newValues[$r - 1]\n    $colCount = $rowValues.Count\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $table.Cell($r, $c)\n        $cell.Range.Text = $rowValues[$c - 1]\n    }\n}\n"}
